# BalaRaju - Completed TDS
# Update TDS (U), total_deducations (W) and NetPay (X) values for rows 2-4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("U2").Value = 1530.15
$ws.Range("W2").Value = 3675.27
$ws.Range("X2").Value = 27474.7

$ws.Range("U3").Value = 2703.12
$ws.Range("W3").Value = 5920.81
$ws.Range("X3").Value = 40804.2

$ws.Range("U4").Value = 13830.8
$ws.Range("W4").Value = 19193.6
$ws.Range("X4").Value = 58681.4
